$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 86.29678392075563
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 94.99655889021044

# Row 3
$ws.Range("B3").Value = 0.06328177979961902
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 16.98373111632243
$ws.Range("E3").Value = 246.9852506941017
$ws.Range("G3").Value = 264.3698484262322

# Row 4
$ws.Range("B4").Value = 0.00006486019690155054
$ws.Range("C4").Value = 0.004309184025731883
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 7.198534972508243

# Row 5
$ws.Range("B5").Value = 3.182878228561681
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.7127328510149897
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("G5").Value = 12.0302756157461
